$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Send")

# Before: H1=USDT, I1=USDC, J1=USDCe  (dimension A1:J1)
# After:  H1=DAI, I1=USDC, J1=USDCe, K1=USDT, L1=WETH  (dimension A1:L1)

# Set the new header values in order.
$ws.Range("H1").Value = "DAI"
$ws.Range("I1").Value = "USDC"
$ws.Range("J1").Value = "USDCe"
$ws.Range("K1").Value = "USDT"
$ws.Range("L1").Value = "WETH"

# Make sure the new cells (K1, L1) carry the same header formatting as the
# rest of the row (bold font, centered, thin box border) - matching style "s=1".
$headerRange = $ws.Range("K1:L1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous
$headerRange.Borders.Weight = 2            # xlThin
